$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.002" or "30.668.33" are not converted to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.668.33'
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").Value = '1.879.19'
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = '239.09'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").Value = '0.4811'
$ws.Range("E7").Value = '  -0.63%  '

$ws.Range("D8").Value = '0.2835'
$ws.Range("E8").Value = '  -2.08%  '

$ws.Range("D9").Value = '0.06524'
$ws.Range("E9").Value = '  -1.24%  '

$ws.Range("D10").Value = '1.953.99'
$ws.Range("E10").Value = '  +3.35%  '

$ws.Range("D11").Value = '0.07491'
$ws.Range("E11").Value = '  +1.13%  '

$ws.Range("D12").Value = '16.55'
$ws.Range("E12").Value = '  -1.50%  '

$ws.Range("D13").Value = '5.099'
$ws.Range("E13").Value = '  -1.67%  '

$ws.Range("D14").Value = '88.25'
$ws.Range("E14").Value = '  -0.39%  '

$ws.Range("D15").Value = '0.6640'
$ws.Range("E15").Value = '  +0.37%  '

$ws.Range("D16").Value = '30.623.08'
$ws.Range("E16").Value = '  +0.63%  '

$ws.Range("D17").Value = '13.30'
$ws.Range("E17").Value = '  -1.74%  '

$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").Value = '2.198.96'
$ws.Range("E19").Value = '  +2.84%  '

$ws.Range("D20").Value = '0.000007602'
$ws.Range("E20").Value = '  -2.17%  '

$ws.Range("D21").Value = '230.11'
$ws.Range("E21").Value = '  +3.97%  '

$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.282'
$ws.Range("E23").Value = '  -1.57%  '

$ws.Range("D24").Value = '6.162'
$ws.Range("E24").Value = '  -0.92%  '

$ws.Range("D25").Value = '168.71'
$ws.Range("E25").Value = '  +3.51%  '

$ws.Range("D26").Value = '9.292'
$ws.Range("E26").Value = '  -0.90%  '

$ws.Range("D27").Value = '18.59'
$ws.Range("E27").Value = '  -1.55%  '

$ws.Range("D28").Value = '1.935'
$ws.Range("E28").Value = '  -0.40%  '

$ws.Range("D29").Value = '1.407'
$ws.Range("E29").Value = '  -2.98%  '

$ws.Range("D30").Value = '0.09727'
$ws.Range("E30").Value = '  +5.56%  '

$ws.Range("D31").Value = '4.340'
$ws.Range("E31").Value = '  +0.23%  '

$ws.Range("D32").Value = '4.011'
$ws.Range("E32").Value = '  -0.68%  '

$ws.Range("E33").Value = '  +0.26%  '

$ws.Range("D34").Value = '1.219'
$ws.Range("E34").Value = '  +5.87%  '

$ws.Range("D35").Value = '0.7502'
$ws.Range("E35").Value = '  -1.21%  '

$ws.Range("D36").Value = '2.716'
$ws.Range("E36").Value = '  +0.48%  '

$ws.Range("D37").Value = '0.01864'
$ws.Range("E37").Value = '  -0.82%  '

$ws.Range("D38").Value = '2.640'
$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '0.9137'
$ws.Range("E39").Value = '  -0.64%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = '2.078'
$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("D41").Value = '106.37'
$ws.Range("E41").Value = '  +0.20%  '

$ws.Range("D42").Value = '0.4275'
$ws.Range("E42").Value = '  -1.86%  '

$ws.Range("D43").Value = '5.776'
$ws.Range("E43").Value = '  -3.36%  '

$ws.Range("D44").Value = '0.9996'
$ws.Range("E44").Value = '  -0.38%  '

$ws.Range("D45").Value = '7.347'
$ws.Range("E45").Value = '  -3.87%  '

$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1291'
$ws.Range("E46").Value = '  -2.62%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '64.34'
$ws.Range("E47").Value = '  -1.83%  '

$ws.Range("D48").Value = '1.477'
$ws.Range("E48").Value = '  -7.92%  '

$ws.Range("D49").Value = '8.955'
$ws.Range("E49").Value = '  +0.04%  '

$ws.Range("D50").Value = '33.79'
$ws.Range("E50").Value = '  -2.25%  '

$ws.Range("D51").Value = '0.05655'
$ws.Range("E51").Value = '  -1.03%  '
